$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

$ws.Cells.Item($row, 4).Value = 44656
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101007
$ws.Cells.Item($row, 10).Value = "Kiwi"
$ws.Cells.Item($row, 11).Value = "Hayward"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 270
$ws.Cells.Item($row, 14).Value = 19000
$ws.Cells.Item($row, 15).Value = 20000
$ws.Cells.Item($row, 16).Value = 19500
$ws.Cells.Item($row, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1083
$ws.Cells.Item($row, 20).Value = 18
